# Update odds for row 2 (txT1cqOG - Sydney FC vs Macarthur FC) - values only change,
# the fixture itself stays in row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 4.2
$ws.Range("L2").Value = 4.33
$ws.Range("N2").Value = 19
$ws.Range("W2").Value = 11
$ws.Range("AB2").Value = 19
$ws.Range("AJ2").Value = 15
$ws.Range("AM2").Value = 29
$ws.Range("AO2").Value = 8.5
$ws.Range("AQ2").Value = 23
$ws.Range("AW2").Value = 301

# The JAPAN - J1 LEAGUE fixture (Urawa Reds vs Sanfrecce Hiroshima, row 3) is removed
# from this week's report entirely; deleting the whole row shifts every later row
# (and the sheet's used range / dimension) up by one automatically.
$ws.Rows(3).Delete()

# After the shift, the South Korea - K League 1 fixtures that used to sit in rows
# 4 and 5 are now in rows 3 and 4 - refresh their updated odds.
$ws.Range("G3").Value = 2.7
$ws.Range("I3").Value = 2.75
$ws.Range("J3").Value = 3.5
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 3.5
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 7.5
$ws.Range("X3").Value = 12
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 29
$ws.Range("AI3").Value = 12
$ws.Range("AO3").Value = 17
$ws.Range("AT3").Value = 2.5
$ws.Range("AX3").Value = 4.5

$ws.Range("G4").Value = 3.4
$ws.Range("I4").Value = 2.15
$ws.Range("J4").Value = 4
$ws.Range("N4").Value = 10
$ws.Range("Q4").Value = 2.05
$ws.Range("R4").Value = 1.75
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.91
$ws.Range("AA4").Value = 29
$ws.Range("AE4").Value = 15
$ws.Range("AG4").Value = 251
$ws.Range("AH4").Value = 7.5
$ws.Range("AI4").Value = 10
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 67
$ws.Range("AT4").Value = 2.63
$ws.Range("AZ4").Value = 23
